$wb = $excel.ActiveWorkbook

# Activate the "Data" worksheet (it becomes the tab-selected / active sheet)
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

# Copy the values from column B (MockBoard) into column C (Board) for rows 2-401
$ws.Range("C2:C401").Value = $ws.Range("B2:B401").Value()

# Select C2:C401 with C2 as the active cell, matching the resulting selection
$ws.Range("C2:C401").Select()
